{"js": "const replacements = [\n  [\"2023-04-12 Wednesday\", \"2023-04-13 Thursday\"],\n  [\"4+37=41\", \"77-67=10\"],\n  [\"32-27=5\", \"23+36=59\"],\n  [\"75-30=45\", \"59+24=83\"],\n  [\"44+39=83\", \"68+28=96\"],\n  [\"32+42=74\", \"36+7=43\"],\n  [\"98-83=15\", \"32+12=44\"],\n  [\"88-49=39\", \"0+81=81\"],\n  [\"83-3=80\", \"67-55=12\"],\n  [\"59-27=32\", \"43-7=36\"],\n  [\"97-22=75\", \"4+93=97\"],\n  [\"53-32=21\", \"23+11=34\"],\n  [\"11+67=78\", \"61+28=89\"],\n  [\"10+81=91\", \"33+39=72\"],\n  [\"75-31=44\", \"24+31=55\"],\n  [\"28+29=57\", \"63+13=76\"],\n  [\"99-88=11\", \"60-30=30\"],\n  [\"93-69=24\", \"87-23=64\"],\n  [\"46+36=82\", \"23+48=71\"],\n  [\"37+61=98\", \"80-71=9\"],\n  [\"60-41=19\", \"6+36=42\"],\n  [\"77-54=23\", \"6+12=18\"],\n  [\"40-29=11\", \"86-62=24\"],\n  [\"89-5=84\", \"2+68=70\"],\n  [\"85-53=32\", \"37+57=94\"],\n  [\"96-71=25\", \"59+20=79\"],\n  [\"96-59=37\", \"54-50=4\"],\n  [\"10+63=73\", \"44+50=94\"],\n  [\"87-44=43\", \"41+3=44\"],\n  [\"66+10=76\", \"44+46=90\"],\n  [\"6+34=40\", \"68+26=94\"],\n  [\"12+34=46\", \"9+6=15\"],\n  [\"60+23=83\", \"60-0=60\"],\n  [\"23+68=91\", \"48-32=16\"],\n  [\"28+27=55\", \"39-28=11\"],\n  [\"51+19=70\", \"87-48=39\"],\n  [\"58+9=67\", \"54+8=62\"],\n  [\"4+47=51\", \"20-11=9\"],\n  [\"46+14=60\", \"41+22=63\"],\n  [\"15-11=4\", \"54+11=65\"],\n  [\"34-5=29\", \"50-9=41\"],\n  [\"42+11=53\", \"7+83=90\"],\n  [\"10+61=71\", \"99-90=9\"],\n  [\"98-76=22\", \"20+2=22\"],\n  [\"95-78=17\", \"32+22=54\"],\n  [\"38-19=19\", \"95-55=40\"],\n  [\"68-60=8\", \"7+26=33\"],\n  [\"81-64=17\", \"34-6=28\"],\n  [\"26+9=35\", \"26+12=38\"],\n  [\"71-49=22\", \"90-7=83\"],\n  [\"17+26=43\", \"11+9=20\"],\n  [\"44+14=58\", \"97-31=66\"],\n  [\"82-3=79\", \"57+41=98\"],\n  [\"53+38=91\", \"59+4=63\"],\n  [\"27-27=0\", \"2+78=80\"],\n  [\"45+51=96\", \"47+35=82\"],\n  [\"31+4=35\", \"42+18=60\"],\n  [\"48-26=22\", \"54+10=64\"],\n  [\"62-13=49\", \"70-54=16\"],\n  [\"67-16=51\", \"93-55=38\"],\n  [\"10-3=7\", \"57-12=45\"],\n  [\"51+33=84\", \"16+59=75\"],\n  [\"91-34=57\", \"70+10=80\"],\n  [\"30+69=99\", \"88-7=81\"],\n  [\"48+3=51\", \"54-21=33\"],\n  [\"73-41=32\", \"62+23=85\"],\n  [\"6+86=92\", \"35+41=76\"],\n  [\"6+48=54\", \"1+17=18\"],\n  [\"79-78=1\", \"64-28=36\"],\n  [\"54+29=83\", \"74+8=82\"],\n  [\"25+14=39\", \"35+52=87\"],\n  [\"7+70=77\", \"7+34=41\"],\n  [\"55-22=33\", \"42-37=5\"],\n  [\"27+37=64\", \"24+3=27\"],\n  [\"70-18=52\", \"87-77=10\"],\n  [\"87-76=11\", \"27+38=65\"],\n  [\"70+22=92\", \"64+10=74\"],\n  [\"93-51=42\", \"37-12=25\"],\n  [\"24+38=62\", \"23+11=34\"],\n  [\"48-13=35\", \"97-54=43\"],\n  [\"10+33=43\", \"21+11=32\"],\n  [\"87-3=84\", \"9+44=53\"],\n  [\"92-90=2\", \"94-67=27\"],\n  [\"24+4=28\", \"43+36=79\"],\n  [\"23-1=22\", \"31+46=77\"],\n  [\"88-60=28\", \"52+18=70\"],\n  [\"64-18=46\", \"14+26=40\"],\n  [\"80+5=85\", \"6+93=99\"],\n  [\"9+16=25\", \"19+80=99\"],\n  [\"90-53=37\", \"24+17=41\"],\n  [\"98-90=8\", \"2+30=32\"],\n  [\"88-17=71\", \"41-22=19\"],\n  [\"14+17=31\", \"15+31=46\"],\n  [\"0+83=83\", \"51-31=20\"],\n  [\"29+36=65\", \"73-69=4\"],\n  [\"79-48=31\", \"24-17=7\"],\n  [\"63-9=54\", \"34-17=17\"],\n  [\"97-66=31\", \"34+61=95\"],\n  [\"92-60=32\", \"97-38=59\"],\n  [\"23+26=49\", \"72+1=73\"],\n  [\"2+72=74\", \"52-1=51\"],\n];\n\nconst body = context.document.body;\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + before);\n  }\n  results.items[0].insertText(after, Word.InsertLocation.replace);\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n$replacements = @(\n    @('2023-04-12 Wednesday', '2023-04-13 Thursday'),\n    @('4+37=41', '77-67=10'),\n    @('32-27=5', '23+36=59'),\n    @('75-30=45', '59+24=83'),\n    @('44+39=83', '68+28=96'),\n    @('32+42=74', '36+7=43'),\n    @('98-83=15', '32+12=44'),\n    @('88-49=39', '0+81=81'),\n    @('83-3=80', '67-55=12'),\n    @('59-27=32', '43-7=36'),\n    @('97-22=75', '4+93=97'),\n    @('53-32=21', '23+11=34'),\n    @('11+67=78', '61+28=89'),\n    @('10+81=91', '33+39=72'),\n    @('75-31=44', '24+31=55'),\n    @('28+29=57', '63+13=76'),\n    @('99-88=11', '60-30=30'),\n    @('93-69=24', '87-23=64'),\n    @('46+36=82', '23+48=71'),\n    @('37+61=98', '80-71=9'),\n    @('60-41=19', '6+36=42'),\n    @('77-54=23', '6+12=18'),\n    @('40-29=11', '86-62=24'),\n    @('89-5=84', '2+68=70'),\n    @('85-53=32', '37+57=94'),\n    @('96-71=25', '59+20=79'),\n    @('96-59=37', '54-50=4'),\n    @('10+63=73', '44+50=94'),\n    @('87-44=43', '41+3=44'),\n    @('66+10=76', '44+46=90'),\n    @('6+34=40', '68+26=94'),\n    @('12+34=46', '9+6=15'),\n    @('60+23=83', '60-0=60'),\n    @('23+68=91', '48-32=16'),\n    @('28+27=55', '39-28=11'),\n    @('51+19=70', '87-48=39'),\n    @('58+9=67', '54+8=62'),\n    @('4+47=51', '20-11=9'),\n    @('46+14=60', '41+22=63'),\n    @('15-11=4', '54+11=65'),\n    @('34-5=29', '50-9=41'),\n    @('42+11=53', '7+83=90'),\n    @('10+61=71', '99-90=9'),\n    @('98-76=22', '20+2=22'),\n    @('95-78=17', '32+22=54'),\n    @('38-19=19', '95-55=40'),\n    @('68-60=8', '7+26=33'),\n    @('81-64=17', '34-6=28'),\n    @('26+9=35', '26+12=38'),\n    @('71-49=22', '90-7=83'),\n    @('17+26=43', '11+9=20'),\n    @('44+14=58', '97-31=66'),\n    @('82-3=79', '57+41=98'),\n    @('53+38=91', '59+4=63'),\n    @('27-27=0', '2+78=80'),\n    @('45+51=96', '47+35=82'),\n    @('31+4=35', '42+18=60'),\n    @('48-26=22', '54+10=64'),\n    @('62-13=49', '70-54=16'),\n    @('67-16=51', '93-55=38'),\n    @('10-3=7', '57-12=45'),\n    @('51+33=84', '16+59=75'),\n    @('91-34=57', '70+10=80'),\n    @('30+69=99', '88-7=81'),\n    @('48+3=51', '54-21=33'),\n    @('73-41=32', '62+23=85'),\n    @('6+86=92', '35+41=76'),\n    @('6+48=54', '1+17=18'),\n    @('79-78=1', '64-28=36'),\n    @('54+29=83', '74+8=82'),\n    @('25+14=39', '35+52=87'),\n    @('7+70=77', '7+34=41'),\n    @('55-22=33', '42-37=5'),\n    @('27+37=64', '24+3=27'),\n    @('70-18=52', '87-77=10'),\n    @('87-76=11', '27+38=65'),\n    @('70+22=92', '64+10=74'),\n    @('93-51=42', '37-12=25'),\n    @('24+38=62', '23+11=34'),\n    @('48-13=35', '97-54=43'),\n    @('10+33=43', '21+11=32'),\n    @('87-3=84', '9+44=53'),\n    @('92-90=2', '94-67=27'),\n    @('24+4=28', '43+36=79'),\n    @('23-1=22', '31+46=77'),\n    @('88-60=28', '52+18=70'),\n    @('64-18=46', '14+26=40'),\n    @('80+5=85', '6+93=99'),\n    @('9+16=25', '19+80=99'),\n    @('90-53=37', '24+17=41'),\n    @('98-90=8', '2+30=32'),\n    @('88-17=71', '41-22=19'),\n    @('14+17=31', '15+31=46'),\n    @('0+83=83', '51-31=20'),\n    @('29+36=65', '73-69=4'),\n    @('79-48=31', '24-17=7'),\n    @('63-9=54', '34-17=17'),\n    @('97-66=31', '34+61=95'),\n    @('92-60=32', '97-38=59'),\n    @('23+26=49', '72+1=73'),\n    @('2+72=74', '52-1=51'),\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $result) {\n        throw \"Find/Replace failed for: $findText\"\n    }\n}"}
